# Generate Report for Handback
# Refresh handback-status timestamps (and the "ht" -> "mt" status code) for
# the c09cf10b.../d7454abb... entries to reflect a newly generated report.
# Because the underlying XLSX stores these repeated values as shared
# strings, every cell that held the old text must be moved to the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" (Overview col G) and the corresponding
# "Correspond Handoff Datetime" on the de-de sheet (col H) all shared the
# string "2016-08-17 02:14:23" -> "2016-08-17 02:15:13".
$wsOverview.Range("G3").Value = "2016-08-17 02:15:13"
$wsOverview.Range("G4").Value = "2016-08-17 02:15:13"
$wsDeDe.Range("H3").Value = "2016-08-17 02:15:13"
$wsDeDe.Range("H4").Value = "2016-08-17 02:15:13"

# Status column ("Priority"/ht-mt code) on both language sheets shared the
# string "ht" -> "mt".
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn "Correspond Handoff Datetime" (col H): "2016-08-17 02:14:18" -> "2016-08-17 02:15:03"
$wsZhCn.Range("H3").Value = "2016-08-17 02:15:03"
$wsZhCn.Range("H4").Value = "2016-08-17 02:15:03"

# zh-cn "Correspond Handback DateTime" (col K): "2016-08-17 02:14:34" -> "2016-08-17 02:15:30"
$wsZhCn.Range("K3").Value = "2016-08-17 02:15:30"
$wsZhCn.Range("K4").Value = "2016-08-17 02:15:30"

# de-de "Correspond Handback DateTime" (col K): "2016-08-17 02:14:41" -> "2016-08-17 02:15:37"
$wsDeDe.Range("K3").Value = "2016-08-17 02:15:37"
$wsDeDe.Range("K4").Value = "2016-08-17 02:15:37"
